# edit.ps1 - applies the resume diff via Word COM-interop (PowerShell style)
#
# Strategy: locate each edited span with Range.Find scoped to the owning
# paragraph (avoids cross-paragraph ambiguity), replace its text, then
# (when the target span must become several runs) re-split the freshly
# written text into sibling runs at the desired offsets.
#
# Word normally coalesces adjacent runs that end up with identical
# formatting after a text edit. Toggling a boolean character property
# (Italic) to a different value and straight back to its original value
# around a text write breaks that coalescing without leaving any visible
# / serialized trace, because Italic=0 is already explicit in every rPr
# in this document. We use that trick both to keep a freshly written run
# isolated from whatever follows it in the paragraph, and to cut a
# longer run into several shorter sibling runs.

$d = $word.ActiveDocument

function Replace-Span {
    param($Para, $OldText, $NewParts)

    $newText = [string]::Join("", $NewParts)

    $rng = $Para.Range
    $found = $rng.Find.Execute($OldText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Host "NOT FOUND: [$OldText]"
        return
    }

    # Write the new text, isolating the resulting run from whatever
    # comes after it in the paragraph (prevents auto-merge).
    $rng.Italic = 1
    $rng.Text = $newText
    $rng.Italic = 0

    # Now cut the written range into sibling runs matching $NewParts,
    # each boundary gets the same isolate trick.
    $base = $rng.Start
    $offset = 0
    for ($i = 0; $i -lt $NewParts.Length - 1; $i++) {
        $offset = $offset + $NewParts[$i].Length
        $piece = $d.Range($base, $base + $offset)
        $piece.Italic = 1
        $piece.Italic = 0
    }
}

# ---------------------------------------------------------------------
# Group A - "Technologies:" line (Skills section)
# ---------------------------------------------------------------------
$pTech = $d.Paragraphs.Item(12)
Replace-Span $pTech "React.js, Node.js, Express.js, " @("React, ", "Node.js,", " ")
Replace-Span $pTech "Redux, " @("Sass, ", "Redux", ", ")
Replace-Span $pTech ", S3" @("")

# ---------------------------------------------------------------------
# Group B - "Led frontend ... Node.js" bullet (WorkHound bullet 1)
# ---------------------------------------------------------------------
$pLed = $d.Paragraphs.Item(17)
Replace-Span $pLed "Led frontend implementation of a new analytics dashboard and developed APIs using Next.js, React.js, " @("Led front", "end ", "development", " of a new analytics dashboard ", "and ", "built", " APIs ", "using Next.js, ", "React, ")

# ---------------------------------------------------------------------
# Group C - "Action that automatically ... notifies slack for" bullet
# ---------------------------------------------------------------------
$pSlack = $d.Paragraphs.Item(18)
Replace-Span $pSlack " Action that automatically creates a ticket in Jira and notifies slack for " @(" Action that automatically creates a ticket in Jira and notifies ", "S", "lack for ")

# ---------------------------------------------------------------------
# Group D - "to deliver features in a fast-paced" bullet
# ---------------------------------------------------------------------
$pDeliver = $d.Paragraphs.Item(19)
Replace-Span $pDeliver "to deliver features in a fast-paced " @("to deliver ", "features in a fast-paced ")

# ---------------------------------------------------------------------
# Group E - "that power the analytics dashboard" bullet
# ---------------------------------------------------------------------
$pDashboard = $d.Paragraphs.Item(20)
Replace-Span $pDashboard "that power the analytics dashboard" @("for", " ", "customer", "-", "facing ", "dashboard", "s")

Write-Host "All edits applied."
